# Auto-generated Excel COM-interop edit script
# Applies cell value updates per the authoritative diff of Leviathan_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3440.2
$ws.Range("I64").Value = 2800.25
$ws.Range("K64").Value = 2800.25
$ws.Range("M64").Value = -2552.25

$ws.Range("H67").Value = 3440.2
$ws.Range("I67").Value = 2800.25
$ws.Range("K67").Value = 2800.25
$ws.Range("M67").Value = -1942.25

$ws.Range("H74").Value = 22004
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 22004
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H107").Value = 54163.285
$ws.Range("I107").Value = 482.16666
$ws.Range("K107").Value = 482.16666
$ws.Range("M107").Value = 1437.83334

$ws.Range("H116").Value = 4227.5
$ws.Range("J116").Value = 4680
$ws.Range("L116").Value = 4680
$ws.Range("N116").Value = -11564

$ws.Range("H132").Value = 1959.9231
$ws.Range("I132").Value = 923.54
$ws.Range("K132").Value = 2770.62
$ws.Range("M132").Value = -240.6199999999999

$ws.Range("H141").Value = 2400
$ws.Range("I141").Value = 2400
$ws.Range("K141").Value = 7200
$ws.Range("M141").Value = -2020

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2506.6875
$ws.Range("I2").Value = 1759.9166
$ws.Range("K2").Value = 1759.9166
$ws.Range("M2").Value = -1646.9166

$ws.Range("H45").Value = 10498.143
$ws.Range("I45").Value = 16558
$ws.Range("K45").Value = 16558
$ws.Range("M45").Value = -16181

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H74").Value = 2567.0625
$ws.Range("I74").Value = 2343
$ws.Range("J74").Value = 3060
$ws.Range("K74").Value = 2343
$ws.Range("L74").Value = 3060
$ws.Range("M74").Value = -1469
$ws.Range("N74").Value = -4808

$ws.Range("H77").Value = 2567.0625
$ws.Range("I77").Value = 2343
$ws.Range("J77").Value = 3060
$ws.Range("K77").Value = 11715
$ws.Range("L77").Value = 15300
$ws.Range("M77").Value = -7347
$ws.Range("N77").Value = -24036

$ws.Range("H116").Value = 2506.6875
$ws.Range("I116").Value = 1759.9166
$ws.Range("K116").Value = 1759.9166
$ws.Range("M116").Value = 534.0834

$ws.Range("H132").Value = 7352.237
$ws.Range("I132").Value = 7639.8286
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 22919.4858
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -20389.4858
$ws.Range("N132").Value = -17051

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2506.6875
$ws.Range("I3").Value = 1759.9166
$ws.Range("K3").Value = 1759.9166
$ws.Range("M3").Value = -1645.9166

$ws.Range("H86").Value = 6430.1
$ws.Range("I86").Value = 5537.75
$ws.Range("J86").Value = 9999.5
$ws.Range("K86").Value = 5537.75
$ws.Range("L86").Value = 9999.5
$ws.Range("M86").Value = -4414.75
$ws.Range("N86").Value = -12245.5

$ws.Range("H89").Value = 6430.1
$ws.Range("I89").Value = 5537.75
$ws.Range("J89").Value = 9999.5
$ws.Range("K89").Value = 27688.75
$ws.Range("L89").Value = 49997.5
$ws.Range("M89").Value = -22072.75
$ws.Range("N89").Value = -61229.5

$ws.Range("H105").Value = 1859.8889
$ws.Range("I105").Value = 1861.25
$ws.Range("K105").Value = 1861.25
$ws.Range("M105").Value = -114.25

$ws.Range("H134").Value = 95063.91
$ws.Range("I134").Value = 116626.766
$ws.Range("K134").Value = 349880.298
$ws.Range("M134").Value = -347345.298

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3930.2407
$ws.Range("I31").Value = 2562.75
$ws.Range("K31").Value = 2562.75
$ws.Range("M31").Value = -2267.75

$ws.Range("H34").Value = 3930.2407
$ws.Range("I34").Value = 2562.75
$ws.Range("K34").Value = 2562.75
$ws.Range("M34").Value = -2360.75

$ws.Range("H58").Value = 1453.8572
$ws.Range("I58").Value = 1451.55
$ws.Range("K58").Value = 1451.55
$ws.Range("M58").Value = -1248.55

$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26498

$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -82488

$ws.Range("H74").Value = 58562.5
$ws.Range("J74").Value = 58275
$ws.Range("L74").Value = 58275
$ws.Range("N74").Value = -60023

$ws.Range("H77").Value = 58562.5
$ws.Range("J77").Value = 58275
$ws.Range("L77").Value = 174825
$ws.Range("N77").Value = -183561

$ws.Range("H132").Value = 5196.1304
$ws.Range("I132").Value = 5319.5293
$ws.Range("K132").Value = 15958.5879
$ws.Range("M132").Value = -13428.5879

$ws.Range("H136").Value = 1453.8572
$ws.Range("I136").Value = 1451.55
$ws.Range("K136").Value = 4354.65
$ws.Range("M136").Value = -1804.65

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H22").Value = 3390
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3390
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 10170
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -10508

$ws.Range("H27").Value = 3390
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3390
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 10170
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -10374

$ws.Range("H32").Value = 687884.6
$ws.Range("I32").Value = 168016.33
$ws.Range("J32").Value = 1311726.6
$ws.Range("K32").Value = 504048.99
$ws.Range("L32").Value = 3935179.8
$ws.Range("M32").Value = -503765.99
$ws.Range("N32").Value = -3935745.8

$ws.Range("H80").Value = 5375
$ws.Range("I80").Value = 3500
$ws.Range("K80").Value = 10500
$ws.Range("M80").Value = -9564

$ws.Range("H83").Value = 5375
$ws.Range("I83").Value = 3500
$ws.Range("K83").Value = 31500
$ws.Range("M83").Value = -26820

$ws.Range("H100").Value = 7423.4614
$ws.Range("I100").Value = 3404
$ws.Range("K100").Value = 10212
$ws.Range("M100").Value = -9401

$ws.Range("H101").Value = 18112.111
$ws.Range("J101").Value = 18112.111
$ws.Range("L101").Value = 54336.333
$ws.Range("N101").Value = -59204.333

$ws.Range("H105").Value = 14968.4
$ws.Range("J105").Value = 14968.4
$ws.Range("L105").Value = 44905.2
$ws.Range("N105").Value = -50147.2

$ws.Range("H121").Value = 18631702
$ws.Range("I121").Value = 66866836
$ws.Range("J121").Value = 79726
$ws.Range("K121").Value = 200600508
$ws.Range("L121").Value = 239178
$ws.Range("M121").Value = -200599198
$ws.Range("N121").Value = -241798

$ws.Range("H129").Value = 98485.48
$ws.Range("J129").Value = 3924.3333
$ws.Range("L129").Value = 11772.9999
$ws.Range("N129").Value = -21772.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 734.6429000000001
$ws.Range("J2").Value = 1173.875
$ws.Range("L2").Value = 1173.875
$ws.Range("N2").Value = -1399.875

$ws.Range("H132").Value = 6980.231
$ws.Range("I132").Value = 4313.684
$ws.Range("J132").Value = 14218
$ws.Range("K132").Value = 12941.052
$ws.Range("L132").Value = 42654
$ws.Range("M132").Value = -10411.052
$ws.Range("N132").Value = -47714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 23051.715
$ws.Range("I46").Value = 37799.082
$ws.Range("K46").Value = 37799.082
$ws.Range("M46").Value = -37611.082

$ws.Range("H93").Value = 14772.434
$ws.Range("I93").Value = 3492.077
$ws.Range("K93").Value = 3492.077
$ws.Range("M93").Value = -2244.077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20499.166
$ws.Range("J4").Value = 20600
$ws.Range("L4").Value = 20600
$ws.Range("N4").Value = -20826

$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960

$ws.Range("H136").Value = 1371.5333
$ws.Range("I136").Value = 1331.6364
$ws.Range("J136").Value = 1481.25
$ws.Range("K136").Value = 3994.9092
$ws.Range("L136").Value = 4443.75
$ws.Range("M136").Value = -1444.9092
$ws.Range("N136").Value = -9543.75
